$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44973
$ws.Range("M2").Value = 55
$ws.Range("N2").Value = 28000
$ws.Range("O2").Value = 28000
$ws.Range("P2").Value = 28000
$ws.Range("S2").Value = 3500

# Row 3
$ws.Range("D3").Value = 44980
$ws.Range("M3").Value = 50
$ws.Range("N3").Value = 25000
$ws.Range("O3").Value = 25000
$ws.Range("P3").Value = 25000
$ws.Range("S3").Value = 3125

# Row 4
$ws.Range("D4").Value = 44971
$ws.Range("M4").Value = 25

# Row 5
$ws.Range("D5").Value = 44981
$ws.Range("M5").Value = 30
